$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sayfa1")
$ws.Range("B6").Value = "süre bitti"
$ws.Range("B6").Select()
